$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.511.74"
$ws.Range("E2").Value = "  -0.30%  "
$ws.Range("D3").Value = "2.287.04"
$ws.Range("E3").Value = "  -0.27%  "
$ws.Range("E4").Value = "  +0.20%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "112.59"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +16.28%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "267.97"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.30%  "
$ws.Range("E7").Value = "  -0.41%  "
$ws.Range("E8").Value = "  +0.04%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.610"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.25%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "47.31"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +3.22%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0934"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.37%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "8.44"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +8.12%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "15.50"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.71%  "
$ws.Range("D15").Value = "2.631.55"
$ws.Range("E15").Value = "  -0.28%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.837"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.70%  "
$ws.Range("D17").Value = "2.282.97"
$ws.Range("E17").Value = "  -0.15%  "
$ws.Range("D18").Value = "43.473.67"
$ws.Range("E18").Value = "  -0.29%  "
$ws.Range("E19").Value = "  +1.14%  "
$ws.Range("E20").Value = "  +5.56%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "72.24"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.46%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.48"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +2.65%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "232.31"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.24%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.36"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.75%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.82"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +12.96%  "
$ws.Range("E26").Value = "  -0.06%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.34"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.50%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "42.53"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +5.35%  "
$ws.Range("E29").Value = "  -1.26%  "
$ws.Range("E30").Value = "  -0.07%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "176.09"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.30%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "21.57"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.53%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0919"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +3.65%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.48"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.88%  "
$ws.Range("E35").Value = "  +0.23%  "
$ws.Range("E36").Value = "  +6.99%  "
$ws.Range("E37").Value = "  +0.33%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0352"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.78%  "
$ws.Range("E39").Value = "  +10.93%  "
$ws.Range("E40").Value = "  +3.95%  "
$ws.Range("B41").Value = "Algorand"
$ws.Range("C41").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.239"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.57%  "
$ws.Range("B42").Value = "MultiversX"
$ws.Range("C42").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "72.44"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +13.33%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "13.43"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +9.20%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.41"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +5.03%  "
$ws.Range("E45").Value = "  +0.05%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "5.88"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +12.20%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "8.71"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.16%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "102.64"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +4.48%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0999"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.88%  "
$ws.Range("E50").Value = "  +2.54%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.441"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.78%  "
